# Re-maps the shuffled p_department / field shared-string values back onto
# their intended cells, and restores the active selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column E (p_department) corrections ---
$ws.Range("E2").Value = "Oncology"
$ws.Range("E11").Value = "Physical Therapy"
$ws.Range("E15").Value = "Plastic Surgery"
$ws.Range("E17").Value = "Oncology"
$ws.Range("E19").Value = "Dermatology"
$ws.Range("E20").Value = "Neurology"
$ws.Range("E25").Value = "Dermatology"
$ws.Range("E26").Value = "Plastic Surgery"

# --- Column I (doctor field) corrections ---
$ws.Range("I3").Value = "Neurology"
$ws.Range("I4").Value = "Psychiatry"
$ws.Range("I5").Value = "Pediatrics"
$ws.Range("I6").Value = "Plastic Surgery"
$ws.Range("I7").Value = "Cardiology"
$ws.Range("I8").Value = "Gynaecology"
$ws.Range("I9").Value = "Surgery"
$ws.Range("I10").Value = "Psychiatry"
$ws.Range("I12").Value = "Oncology"
$ws.Range("I13").Value = "Pediatrics"
$ws.Range("I14").Value = "Psychiatry"
$ws.Range("I16").Value = "Oncology"
$ws.Range("I18").Value = "Cardiology"
$ws.Range("I21").Value = "Cardiology"
$ws.Range("I22").Value = "Physical Therapy"
$ws.Range("I23").Value = "Endocrinology"
$ws.Range("I24").Value = "Cardiology"

# --- View state: move the active selection (mirrors the saved cursor position) ---
$ws.Activate()
$ws.Range("I17").Select()
